$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "식당"
$ws.Range("C2").Value = 1297
$ws.Range("B3").Value = "메뉴"
$ws.Range("C3").Value = 662
$ws.Range("B4").Value = "없"
$ws.Range("C4").Value = 555
$ws.Range("B5").Value = "밥"
$ws.Range("C5").Value = 478
$ws.Range("B6").Value = "음식"
$ws.Range("C6").Value = 395
$ws.Range("B7").Value = "식단"
$ws.Range("C7").Value = 355
$ws.Range("B8").Value = "좋"
$ws.Range("C8").Value = 355
$ws.Range("B9").Value = "오늘"
$ws.Range("C9").Value = 352
$ws.Range("B10").Value = "식사"
$ws.Range("C10").Value = 350
$ws.Range("B11").Value = "맛있"
$ws.Range("C11").Value = 349
$ws.Range("B12").Value = "중앙"
$ws.Range("C12").Value = 349
$ws.Range("B13").Value = "배식"
$ws.Range("C13").Value = 344
$ws.Range("B14").Value = "맛"
$ws.Range("C14").Value = 330
$ws.Range("B15").Value = "사람"
$ws.Range("C15").Value = 323
$ws.Range("B16").Value = "반찬"
$ws.Range("C16").Value = 302
$ws.Range("B17").Value = "복지관"
$ws.Range("C17").Value = 291
$ws.Range("B18").Value = "식"
$ws.Range("C18").Value = 282
$ws.Range("B19").Value = "분"
$ws.Range("C19").Value = 270
$ws.Range("B20").Value = "시간"
$ws.Range("C20").Value = 252
$ws.Range("B21").Value = "아침"
$ws.Range("C21").Value = 242
$ws.Range("B22").Value = "많"
$ws.Range("C22").Value = 236
$ws.Range("B23").Value = "생각"
$ws.Range("C23").Value = 225
$ws.Range("B24").Value = "개선"
$ws.Range("C24").Value = 208
$ws.Range("B25").Value = "감사"
$ws.Range("C25").Value = 206
$ws.Range("B26").Value = "글"
$ws.Range("C26").Value = 195
$ws.Range("B27").Value = "국"
$ws.Range("C27").Value = 179
$ws.Range("B28").Value = "라면"
$ws.Range("C28").Value = 176
$ws.Range("B29").Value = "점심"
$ws.Range("C29").Value = 169
$ws.Range("B30").Value = "주세"
$ws.Range("C30").Value = 166
$ws.Range("B31").Value = "이용"
$ws.Range("C31").Value = 150
$ws.Range("B32").Value = "김치"
$ws.Range("C32").Value = 146
$ws.Range("B33").Value = "말"
$ws.Range("C33").Value = 144
$ws.Range("B34").Value = "중식"
$ws.Range("C34").Value = 143
$ws.Range("B35").Value = "그렇"
$ws.Range("C35").Value = 142
$ws.Range("B36").Value = "신경"
$ws.Range("C36").Value = 129
$ws.Range("B37").Value = "청운"
$ws.Range("C37").Value = 129
$ws.Range("B38").Value = "안"
$ws.Range("C38").Value = 129
$ws.Range("B39").Value = "아주머니"
$ws.Range("C39").Value = 129
$ws.Range("B40").Value = "정도"
$ws.Range("C40").Value = 128
$ws.Range("B41").Value = "부탁"
$ws.Range("C41").Value = 127
$ws.Range("B42").Value = "중"
$ws.Range("C42").Value = 122
$ws.Range("B43").Value = "청주"
$ws.Range("C43").Value = 121
$ws.Range("B44").Value = "영양사"
$ws.Range("C44").Value = 120
$ws.Range("B45").Value = "전"
$ws.Range("C45").Value = 114
$ws.Range("B46").Value = "기분"
$ws.Range("C46").Value = 112
$ws.Range("B47").Value = "이렇"
$ws.Range("C47").Value = 111
$ws.Range("B48").Value = "기호"
$ws.Range("C48").Value = 110
$ws.Range("B49").Value = "관련"
$ws.Range("C49").Value = 109
$ws.Range("B50").Value = "앞"
$ws.Range("C50").Value = 107
$ws.Range("B51").Value = "조식"
$ws.Range("C51").Value = 107
$ws.Range("B52").Value = "요즘"
$ws.Range("C52").Value = 105
$ws.Range("B53").Value = "업체"
$ws.Range("C53").Value = 104
$ws.Range("B54").Value = "기숙사"
$ws.Range("C54").Value = 102
$ws.Range("B55").Value = "물"
$ws.Range("C55").Value = 101
$ws.Range("B56").Value = "고담"
$ws.Range("C56").Value = 100
$ws.Range("B57").Value = "주말"
$ws.Range("C57").Value = 99
$ws.Range("B58").Value = "금일"
$ws.Range("C58").Value = 97
$ws.Range("B59").Value = "한식"
$ws.Range("C59").Value = 96
$ws.Range("B60").Value = "야식"
$ws.Range("C60").Value = 96
$ws.Range("B61").Value = "저녁"
$ws.Range("C61").Value = 95
$ws.Range("B62").Value = "그릇"
$ws.Range("C62").Value = 93
$ws.Range("B63").Value = "석식"
$ws.Range("C63").Value = 89
$ws.Range("B64").Value = "건의"
$ws.Range("C64").Value = 89
$ws.Range("B65").Value = "경우"
$ws.Range("C65").Value = 88
$ws.Range("B66").Value = "데"
$ws.Range("C66").Value = 87
$ws.Range("B67").Value = "회사"
$ws.Range("C67").Value = 86
$ws.Range("B68").Value = "일"
$ws.Range("C68").Value = 85
$ws.Range("B69").Value = "문제"
$ws.Range("C69").Value = 84
$ws.Range("B70").Value = "수고"
$ws.Range("C70").Value = 83
$ws.Range("B71").Value = "위생"
$ws.Range("C71").Value = 82
$ws.Range("B72").Value = "직원"
$ws.Range("C72").Value = 82
$ws.Range("B73").Value = "볶음"
$ws.Range("C73").Value = 82
$ws.Range("B74").Value = "후"
$ws.Range("C74").Value = 82
$ws.Range("B75").Value = "안녕"
$ws.Range("C75").Value = 80
$ws.Range("B76").Value = "번"
$ws.Range("C76").Value = 80
$ws.Range("B77").Value = "빵"
$ws.Range("C77").Value = 80
$ws.Range("B78").Value = "고기"
$ws.Range("C78").Value = 79
$ws.Range("B79").Value = "조리"
$ws.Range("C79").Value = 78
$ws.Range("B80").Value = "제공"
$ws.Range("C80").Value = 77
$ws.Range("B81").Value = "친절"
$ws.Range("C81").Value = 77
$ws.Range("B82").Value = "준비"
$ws.Range("C82").Value = 76
$ws.Range("B83").Value = "양"
$ws.Range("C83").Value = 76
$ws.Range("B84").Value = "시"
$ws.Range("C84").Value = 75
$ws.Range("B85").Value = "다르"
$ws.Range("C85").Value = 74
$ws.Range("B86").Value = "하세"
$ws.Range("C86").Value = 72
$ws.Range("B87").Value = "터"
$ws.Range("C87").Value = 71
$ws.Range("B88").Value = "답변"
$ws.Range("C88").Value = 70
$ws.Range("B89").Value = "층"
$ws.Range("C89").Value = 70
$ws.Range("B90").Value = "영양"
$ws.Range("C90").Value = 69
$ws.Range("B91").Value = "관리"
$ws.Range("C91").Value = 69
$ws.Range("B92").Value = "어제"
$ws.Range("C92").Value = 68
$ws.Range("B93").Value = "확인"
$ws.Range("C93").Value = 67
$ws.Range("B94").Value = "어떻"
$ws.Range("C94").Value = 67
$ws.Range("B95").Value = "캠퍼스"
$ws.Range("C95").Value = 67
$ws.Range("B96").Value = "면"
$ws.Range("C96").Value = 67
$ws.Range("B97").Value = "불만"
$ws.Range("C97").Value = 67
$ws.Range("B98").Value = "돈"
$ws.Range("C98").Value = 67
$ws.Range("B99").Value = "코너"
$ws.Range("C99").Value = 66
$ws.Range("B100").Value = "질"
$ws.Range("C100").Value = 65
$ws.Range("B101").Value = "샐러드"
$ws.Range("C101").Value = 64
